{"js": "// Fix a copy/paste error in the \"sala\" entity's foreign-key reference:\n// \"#id_mote -> Sala[NN]\" should read \"#id_mote -> Mote[NN]\" (the mote\n// table, not sala, is what #id_mote references).\n//\n// The string \"Sala[NN]\" also appears later in the document (for the\n// genuine sala foreign keys on the \"regras\" and \"atuador\" lines), so\n// the search is scoped to the first paragraph (\"sala (...)\") to change\n// only the intended occurrence.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst salaParagraph = paragraphs.items[0];\n\nconst matches = salaParagraph.search(\"Sala[NN]\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  matches.items[0].insertText(\"Mote[NN]\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix a copy/paste error in the \"sala\" entity's foreign-key reference:\n# \"#id_mote -> Sala[NN]\" should read \"#id_mote -> Mote[NN]\" (the mote\n# table, not sala, is what #id_mote references).\n#\n# \"Sala[NN]\" also occurs later in the document for the genuine sala\n# foreign keys (on the \"regras\" and \"atuador\" lines), so the\n# Find/Replace is scoped to the first paragraph (\"sala (...)\") so only\n# the intended occurrence changes.\n\n$d = $word.ActiveDocument\n\n$salaParagraphRange = $d.Paragraphs(1).Range\n\n$find = $salaParagraphRange.Find\n$find.Text = \"Sala[NN]\"\n$find.Replacement.Text = \"Mote[NN]\"\n$find.Execute(\n    \"Sala[NN]\",   # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    \"Mote[NN]\",   # ReplaceWith\n    2             # Replace (wdReplaceAll)\n)\n"}
